$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 12 ("Score"): the shared formula in E12:G12 used to stop at G;
# extend the same calculation across the full used range (H12:Z12),
# written separately from E12:G12 so the pre-existing shared formula
# there is left intact, then stamp the H:Z block with E12's format.
# ------------------------------------------------------------------
$ws.Range("H12:Z12").Formula = '=$B$6*H6+$B$7*H7+$B$8*H8+$B$9*H9+(H10-H11)*$B$12+$B$13'
$ws.Range("E12").Copy()
$ws.Range("H12:Z12").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Row 13 ("Probability"): fix the formula to reference $B$15 (the
# correctly-derived pi_mult) instead of the raw $A$17 input, and
# extend the shared formula from E13:G13 out to H13:Z13.
# ------------------------------------------------------------------
$ws.Range("D13").Formula = '=$B$15*EXP(D12)/(1+$B$15*EXP(D12))'
$ws.Range("E13:G13").Formula = '=$B$15*EXP(E12)/(1+$B$15*EXP(E12))'
$ws.Range("H13:Z13").Formula = '=$B$15*EXP(H12)/(1+$B$15*EXP(H12))'
$ws.Range("E13").Copy()
$ws.Range("H13:Z13").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Row 14 (percentage label): extend the shared formula from E14:G14
# out to H14:Z14, matching the L14 format (style index used by the
# rest of that row) across the whole E14:Z14 block.
# ------------------------------------------------------------------
$ws.Range("H14:Z14").Formula = '=IF(OR(ISBLANK(H6), ISBLANK(H7), ISBLANK(H8), ISBLANK(H9), ISBLANK(H10), ISBLANK(H11)),"",CONCAT(ROUND(H13*100,0), "%"))'
$ws.Range("L14").Copy()
$ws.Range("E14:Z14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Selection moves to F12 (view stays scrolled to row 4).
# ------------------------------------------------------------------
$ws.Range("F12").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
